$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-17 (Equipment, Manufacturer, Model, Serial No., ID No., Calibration Date, Certificate No.)
$data = @(
    @("HOT AIR OVEN",       "BINDER",          "RF 53",       "RL15-07514",          "NHP-CT-013",    45825, "CH25052997"),
    @("CENTRIFUGE",         "HETTICH",         "ROTINA 46",   "0000240",             "NHP-C-CT-003",  45824, "CF25052998"),
    @("ELECTRONIC BALANCE", "OHAUS",           "V22PWE3T",    "8338347025",          "NHP-HIS-006",   45825, "EB25052999"),
    @("ELECTRONIC BALANCE", "OHAUS",           "FB3",         "8027090157",          "NHP-HIS-007",   45824, "EB25053000"),
    @("ELECTRONIC BALANCE", "OHAUS",           "V22PWE3T",    "8337110661",          "NHP-HIS-008",   45825, "EB25053001"),
    @("ELECTRONIC BALANCE", "SARTORIUS",       "BSA3202S-CW", "26290320",            "NHP-HIS-009",   45824, "EB25053002"),
    @("HOT AIR OVEN",       "MEMMERT",         "UF30",        "B124.2173",           "NHP-HIS-015",   45826, "CH25053003"),
    @("HOT AIR OVEN",       "MEMMERT",         "UN30",        "B121.0649",           "NHP-HIS-016",   45825, "CH25053004"),
    @("WATER BATH",         "ZEEDO",           "HS1125",      "CJHS2099-32HS1125",   "NHP-HIS-033",   45826, "WB25053005"),
    @("WATER BATH",         "ELECTROTHERMAL",  "MH8517",      "M553340/07",          "NHP-HIS-034",   45825, "WB25053006"),
    @("WATER BATH",         "ELECTROTHERMAL",  "MH8517",      "M577610/08",          "NHP-HIS-035",   45826, "WB25053007"),
    @("WATER BATH",         "LEICA",           "HI1210",      "16634",               "NHP-HIS-036",   45825, "WB25053008"),
    @("REFRIGERATOR",       "MEDICOOL",        "MDC145",      "MDC145-202311046",    "NHP-MOL-003",   45824, "CH25053010"),
    @("FREEZER",            "MEDICOOL",        "MDC115",      "MDC115-202311019",    "NHP-MOL-004",   45824, "CH25053011"),
    @("CENTRIFUGE",         "LABTRON",         "CF8",         "1110801241A006",      "NHP-MOL-006",   45826, "CF25053012"),
    @("CENTRIFUGE",         "UGAIYA",          "TXD3",        "W19800024030108",     "NHP-CT-0048",   45824, "CF25053009")
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

# Delete the last row (row 107), shifting cells up
$ws.Rows.Item(107).Delete()

# Update print area to reflect the new extent
$ws.PageSetup.PrintArea = "`$A`$1:`$H`$106"

# Update the selection to match the target state
$ws.Range("H8").Select()
